$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web")

# New "ToolTips" section appended after the existing "UseTabContent" rows.
# Row 159: section header (column A only)
$ws.Range("A159").Value = "ToolTips"

# Row 160: nav locator
$ws.Range("A160").Value = "toolTipsNav"
$ws.Range("B160").Value = "//div[@class='element-list collapse show']//li[@id='item-6']"
$ws.Range("C160").Value = "By.xpath"

# Row 161: tooltip button locator
$ws.Range("A161").Value = "BtnTooltipElement"
$ws.Range("B161").Value = "//button[@id='toolTipButton']"
$ws.Range("C161").Value = "By.xpath"

# Row 162: tooltip text locator
$ws.Range("A162").Value = "tooltipText"
$ws.Range("B162").Value = "//div[@class='tooltip-inner']"
$ws.Range("C162").Value = "By.xpath"

# Move the active selection to the last edited cell, matching the saved view.
$ws.Range("B162").Select()
